# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Sri Lanka" / "Guinea-Bisau" rows (names swap, data updates) ---
# Row 105 keeps its position but now shows "Guinea-Bisau" with refreshed data
$ws.Range("A105").Value = "Guinea-Bisau"
$ws.Range("B105").Value = 990
$ws.Range("C105").Value = 21
$ws.Range("D105").Value = 26
$ws.Range("E105").Value = 960
$ws.Range("H105").Value = 4

# Row 106 keeps its position but now shows "Sri Lanka" with its (unchanged) prior data
$ws.Range("A106").Value = "Sri Lanka"
$ws.Range("B106").Value = 981
$ws.Range("C106").Value = 24
$ws.Range("D106").Value = 538
$ws.Range("E106").Value = 434
$ws.Range("H106").Value = 9

# --- Update "Datos actualizados" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Mayo de 2020 a las 23:05"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1524042
$ws.Range("C4").Value = 16269
$ws.Range("D4").Value = 342844
$ws.Range("E4").Value = 1090305
$ws.Range("G4").Value = 780
$ws.Range("H4").Value = 90893

# --- Row 14: India ---
$ws.Range("B14").Value = 95698
$ws.Range("C14").Value = 5050
$ws.Range("E14").Value = 55878

# --- Row 15: Peru ---
$ws.Range("D15").Value = 28621
$ws.Range("E15").Value = 61004

# --- Row 17: Canada ---
$ws.Range("B17").Value = 76945
$ws.Range("C17").Value = 1081
$ws.Range("D17").Value = 38474
$ws.Range("E17").Value = 32690

# --- Row 101: Maldivas ---
$ws.Range("B101").Value = 1094
$ws.Range("C101").Value = 16
$ws.Range("E101").Value = 1032

# --- Row 110: Niger ---
$ws.Range("B110").Value = 904
$ws.Range("C110").Value = 15
$ws.Range("D110").Value = 698
$ws.Range("E110").Value = 152
$ws.Range("G110").Value = 3
$ws.Range("H110").Value = 54

# --- Row 159: Mozambique ---
$ws.Range("D159").Value = 44
$ws.Range("E159").Value = 93
